# Applies the author's edit: adds a new example row (row 8) to the
# "一般企業" sheet (first / active sheet) and updates the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new row 8, mirroring row 7 but with a distinct label and the
#     previously-blank N/P values filled in. ---
$ws.Range("A8").Value = "範例"
$ws.Range("B8").Value = "台中市"
$ws.Range("C8").Value = "豐原區"
$ws.Range("D8").Value = "菸草製造業"
$ws.Range("E8").Value = 20000
$ws.Range("F8").Value = 160
$ws.Range("G8").Value = 120
$ws.Range("H8").Value = 100000
$ws.Range("I8").Value = "台北市"
$ws.Range("J8").Value = "士林區"
$ws.Range("K8").Value = "是"
$ws.Range("L8").Value = "否"
$ws.Range("M8").Value = 200
$ws.Range("N8").Value = 75
$ws.Range("O8").Value = 120
$ws.Range("P8").Value = 100

# --- Update the saved cursor/selection position on the sheet. ---
$ws.Range("O17").Select() | Out-Null
